$wb = $excel.ActiveWorkbook

# --- Status text change: "Ready for handoff" -> "In Translation" ---
# This shared string is used as the per-file/per-language Status value on
# the Overview sheet (zh-cn / de-de columns) and on each language sheet's
# "Status" column.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- Column width side-effect ---
# Real Excel re-measured ("best fit") these Status columns after the text
# got shorter ("Ready for handoff" -> "In Translation"), shrinking them
# from ~17.22 to ~13.41 characters. Reproduce that resize explicitly.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
